$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '43.627.81'
Set-TextValue 'E2' '  +0.01%  '
Set-TextValue 'D3' '2.290.69'
Set-TextValue 'E3' '  +0.27%  '
Set-TextValue 'E4' '  +0.24%  '
Set-TextValue 'D5' '114.37'
Set-TextValue 'E5' '  +19.65%  '
Set-TextValue 'E6' '  +0.55%  '
Set-TextValue 'D7' '0.624'
Set-TextValue 'E7' '  +0.18%  '
Set-TextValue 'E8' '  -0.01%  '
Set-TextValue 'D9' '0.620'
Set-TextValue 'E9' '  +2.67%  '
Set-TextValue 'D10' '48.17'
Set-TextValue 'E10' '  +4.91%  '
Set-TextValue 'E11' '  +0.02%  '
Set-TextValue 'D12' '8.73'
Set-TextValue 'E12' '  +11.92%  '
Set-TextValue 'D13' '0.106'
Set-TextValue 'E13' '  +0.69%  '
Set-TextValue 'E14' '  +3.12%  '
Set-TextValue 'D15' '2.633.61'
Set-TextValue 'E15' '  +0.20%  '
Set-TextValue 'E16' '  +0.07%  '
Set-TextValue 'D17' '2.288.58'
Set-TextValue 'E17' '  +0.01%  '
Set-TextValue 'D18' '43.571.83'
Set-TextValue 'E19' '  +2.13%  '
Set-TextValue 'E20' '  +5.18%  '
Set-TextValue 'D21' '72.52'
Set-TextValue 'E21' '  +0.53%  '
Set-TextValue 'E22' '  -0.36%  '
Set-TextValue 'D23' '232.84'
Set-TextValue 'E23' '  +0.05%  '
Set-TextValue 'D24' '9.74'
Set-TextValue 'E24' '  +6.11%  '
Set-TextValue 'E25' '  +12.53%  '
Set-TextValue 'E26' '  -0.03%  '
Set-TextValue 'D27' '11.51'
Set-TextValue 'E27' '  +3.52%  '
Set-TextValue 'D28' '42.14'
Set-TextValue 'E28' '  +3.90%  '
Set-TextValue 'E29' '  -1.98%  '
Set-TextValue 'E30' '  -0.22%  '
Set-TextValue 'D31' '176.67'
Set-TextValue 'E31' '  +0.32%  '
Set-TextValue 'D32' '0.0934'
Set-TextValue 'E32' '  +4.79%  '
Set-TextValue 'E33' '  -0.94%  '
Set-TextValue 'D34' '5.54'
Set-TextValue 'E34' '  +3.48%  '
Set-TextValue 'E35' '  +0.92%  '
Set-TextValue 'D36' '4.72'
Set-TextValue 'E36' '  +8.87%  '
Set-TextValue 'E37' '  -0.09%  '
Set-TextValue 'D38' '0.0356'
Set-TextValue 'E38' '  +0.22%  '
Set-TextValue 'D39' '3.84'
Set-TextValue 'E39' '  +12.19%  '
Set-TextValue 'B40' 'Celestia'
Set-TextValue 'C40' 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue 'D40' '13.87'
Set-TextValue 'E40' '  +12.99%  '
Set-TextValue 'D41' '2.43'
Set-TextValue 'E41' '  +5.29%  '
Set-TextValue 'B42' 'Algorand'
Set-TextValue 'C42' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D42' '0.244'
Set-TextValue 'E42' '  +3.81%  '
Set-TextValue 'D43' '72.58'
Set-TextValue 'E43' '  +11.58%  '
Set-TextValue 'E44' '  +7.41%  '
Set-TextValue 'E45' '  +16.39%  '
Set-TextValue 'E46' '  +0.11%  '
Set-TextValue 'E47' '  -0.51%  '
Set-TextValue 'B48' 'Aave'
Set-TextValue 'C48' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D48' '102.53'
Set-TextValue 'E48' '  +5.67%  '
Set-TextValue 'B49' 'Cronos'
Set-TextValue 'C49' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D49' '0.100'
Set-TextValue 'E49' '  -1.46%  '
Set-TextValue 'D50' '1.22'
Set-TextValue 'E50' '  +2.78%  '
Set-TextValue 'D51' '0.450'
Set-TextValue 'E51' '  +3.12%  '
